# Fruta / hortaliza, semanal
# Insert a new weekly record at row 572 (sheet data rows shift down by one:
# the former row 572 becomes row 573, ..., the former row 649 becomes row 650).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 572-649 down to 573-650 and open up a blank row 572.
$ws.Rows.Item(572).Insert()

# Populate the newly inserted row 572 with the new weekly observation.
$ws.Cells.Item(572, 1).Value  = 6
$ws.Cells.Item(572, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(572, 3).Value  = "Metropolitana"
$ws.Cells.Item(572, 4).Value  = 44491
$ws.Cells.Item(572, 5).Value  = 13
$ws.Cells.Item(572, 6).Value  = 100112024
$ws.Cells.Item(572, 7).Value  = "Choclo"
$ws.Cells.Item(572, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(572, 9).Value  = "Primera"
$ws.Cells.Item(572, 10).Value = 630
$ws.Cells.Item(572, 11).Value = 25000
$ws.Cells.Item(572, 12).Value = 30000
$ws.Cells.Item(572, 13).Value = 26984
$ws.Cells.Item(572, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(572, 15).Value = "Argentina"
$ws.Cells.Item(572, 16).Value = 540
$ws.Cells.Item(572, 17).Value = 50
$ws.Cells.Item(572, 18).Value = "Hortaliza"
